# Criando tabela para o chat
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove yellow highlight fill previously applied to B23:F23
$ws.Range("B23:F23").Style = "Normal"

# New "Mensagem" table starting at row 82
$ws.Range("A82").Value = "Mensagem"

$ws.Range("A83").Value = "Chave"
$ws.Range("B83").Value = "Campo"
$ws.Range("C83").Value = "Tipo/Tamanho"
$ws.Range("D83").Value = "Requerido"
$ws.Range("E83").Value = "Regra"
$ws.Range("F83").Value = "Obs"

$ws.Range("A84").Value = "PK"
$ws.Range("B84").Value = "id_Mensagem"
$ws.Range("C84").Value = "int"
$ws.Range("D84").Value = "not null"
$ws.Range("E84").Value = "auto_increment"
$ws.Range("F84").Value = "primary key,"

$ws.Range("A85").Value = "FK"
$ws.Range("B85").Value = "id_Proposta_Mensagem"
$ws.Range("C85").Value = "int"
$ws.Range("D85").Value = "not null"
$ws.Range("F85").Value = ","

$ws.Range("B86").Value = "texto_Mensagem"
$ws.Range("C86").Value = "varchar(255)"
$ws.Range("D86").Value = "not null"
$ws.Range("F86").Value = ","

$ws.Range("B87").Value = "registro_Mensagem"
$ws.Range("C87").Value = "timestamp"
$ws.Range("D87").Value = "not null"
$ws.Range("F87").Value = ","

$ws.Range("B88").Value = "remetente_Mensagem"
$ws.Range("C88").Value = "int"
$ws.Range("D88").Value = "not null"
$ws.Range("F88").Value = ","

# status_Mensagem / obs_Mensagem typed first (rows 89/90)...
$ws.Range("B89").Value = "status_Mensagem"
$ws.Range("C89").Value = "varchar(25)"
$ws.Range("D89").Value = "not null"
$ws.Range("F89").Value = ","

$ws.Range("B90").Value = "obs_Mensagem"
$ws.Range("C90").Value = "varchar(255)"
$ws.Range("D90").Value = "null"

# ... then a row is inserted above them to add arquivo_Mensagem afterwards,
# which is why arquivo_Mensagem ends up with the highest shared-string index.
$ws.Rows("89:89").Insert()

$ws.Range("B89").Value = "arquivo_Mensagem"
$ws.Range("C89").Value = "longblob"
$ws.Range("D89").Value = "null"
$ws.Range("F89").Value = ","

# Update the visible window/selection to match the saved state
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 73
$ws.Range("H84").Select()
